$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "Player Info" worksheet before the first sheet
#    ("ODI Batting") so the final sheet order is:
#      Player Info, ODI Batting, ODI Bowling
# ------------------------------------------------------------------
$battingSheetForInsert = $wb.Worksheets.Item("ODI Batting")
$infoSheet = $wb.Worksheets.Add($battingSheetForInsert)
$infoSheet.Name = "Player Info"

# NOTE: the variable used to position the new sheet (above) becomes a
# stale/positional reference once the new sheet has been inserted in
# front of it, so re-fetch "ODI Batting" fresh, by name, before using
# it any further.
$battingSheet = $wb.Worksheets.Item("ODI Batting")

# Header row
$infoSheet.Range("A1").Value = "ID"
$infoSheet.Range("B1").Value = "NAME"
$infoSheet.Range("C1").Value = "BATTING_HAND"
$infoSheet.Range("D1").Value = "BOWL_STYLE"

# Match the bold / bordered / centered header formatting used by the
# header rows on the other two sheets.
$infoHeader = $infoSheet.Range("A1:D1")
$infoHeader.Font.Bold = $true
$infoHeader.HorizontalAlignment = -4108
$infoHeader.VerticalAlignment = -4160
$infoHeader.Borders.LineStyle = 1

# Make sure the ID column is stored as text (matches original data which
# stores the id as an inline string, not a number) before writing values.
$infoSheet.Range("A2").NumberFormat = "@"

# Data row
$infoSheet.Range("A2").Value = "6556"
$infoSheet.Range("B2").Value = "Christopher M Mcbride"
$infoSheet.Range("C2").Value = "Right Handed"
$infoSheet.Range("D2").Value = "Right Arm Medium"

# ------------------------------------------------------------------
# 2. "ODI Batting" sheet: rename MATCH_CARD_LINK column to MATCH_CODE
#    and replace the full scorecard URL with just the numeric match
#    code extracted from it.
# ------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2:D12").NumberFormat = "@"

$battingCodes = @{
    2  = "4573"
    3  = "4617"
    4  = "4629"
    5  = "4677"
    6  = "4681"
    7  = "4680"
    8  = "4684"
    9  = "4702"
    10 = "4703"
    11 = "4705"
    12 = "4706"
}
foreach ($r in $battingCodes.Keys) {
    $battingSheet.Range("D$r").Value = $battingCodes[$r]
}

# ------------------------------------------------------------------
# 3. "ODI Bowling" sheet: same MATCH_CARD_LINK -> MATCH_CODE change.
# ------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingSheet.Range("B2:B5").NumberFormat = "@"

$bowlingCodes = @{
    2 = "4617"
    3 = "4629"
    4 = "4681"
    5 = "4684"
}
foreach ($r in $bowlingCodes.Keys) {
    $bowlingSheet.Range("B$r").Value = $bowlingCodes[$r]
}
